{"js": "// Map of old cell text -> new cell text, taken from the target diff.\n// Every data cell in the single table of this document is replaced;\n// the mapping is keyed by the *original* text so the script is\n// resilient to re-ordering and only touches cells that actually match.\nconst replacements = {\n  \"68+16=84\": \"95-44=51\",\n  \"81+2=83\": \"14+84=98\",\n  \"54+26=80\": \"60-35=25\",\n  \"48+20=68\": \"4+2=6\",\n  \"24+62=86\": \"20+43=63\",\n  \"11+13=24\": \"27+57=84\",\n  \"7+40=47\": \"25+21=46\",\n  \"77-1=76\": \"81-66=15\",\n  \"59-48=11\": \"67+23=90\",\n  \"99-81=18\": \"20+74=94\",\n  \"4+48=52\": \"37+61=98\",\n  \"78-50=28\": \"30+67=97\",\n  \"56+13=69\": \"8+73=81\",\n  \"70-6=64\": \"94-82=12\",\n  \"70-5=65\": \"12+32=44\",\n  \"24-6=18\": \"35+52=87\",\n  \"37+6=43\": \"31+61=92\",\n  \"60+32=92\": \"32+34=66\",\n  \"34-24=10\": \"29+61=90\",\n  \"10+76=86\": \"58+17=75\",\n  \"62-42=20\": \"5+26=31\",\n  \"30+12=42\": \"68+20=88\",\n  \"32-10=22\": \"33-0=33\",\n  \"42+49=91\": \"34-7=27\",\n  \"23+25=48\": \"62+11=73\",\n  \"75+17=92\": \"37+26=63\",\n  \"41-19=22\": \"26+47=73\",\n  \"44-31=13\": \"63-14=49\",\n  \"76+23=99\": \"17+76=93\",\n  \"57-39=18\": \"93+4=97\",\n  \"58+3=61\": \"99-22=77\",\n  \"68+5=73\": \"42+10=52\",\n  \"7+27=34\": \"25+52=77\",\n  \"21+12=33\": \"43-22=21\",\n  \"32+11=43\": \"16+53=69\",\n  \"97-90=7\": \"71+7=78\",\n  \"44-9=35\": \"15+45=60\",\n  \"97-28=69\": \"49-44=5\",\n  \"99-11=88\": \"69-39=30\",\n  \"67-13=54\": \"57-18=39\",\n  \"23+71=94\": \"92-71=21\",\n  \"5+9=14\": \"55+28=83\",\n  \"38+12=50\": \"54-17=37\",\n  \"82-81=1\": \"66+23=89\",\n  \"58-6=52\": \"74-65=9\",\n  \"14+15=29\": \"42-36=6\",\n  \"2+88=90\": \"96-92=4\",\n  \"28+43=71\": \"56-54=2\",\n  \"69-38=31\": \"75+10=85\",\n  \"8+51=59\": \"11-7=4\",\n  \"93-78=15\": \"3+48=51\",\n  \"79-52=27\": \"35+16=51\",\n  \"86-14=72\": \"21-8=13\",\n  \"10+82=92\": \"96-6=90\",\n  \"58-38=20\": \"60+20=80\",\n  \"72-41=31\": \"55-11=44\",\n  \"57-13=44\": \"68+23=91\",\n  \"45+2=47\": \"63-60=3\",\n  \"76-71=5\": \"57-35=22\",\n  \"49-46=3\": \"61-4=57\",\n  \"17+33=50\": \"86+9=95\",\n  \"19+12=31\": \"65+17=82\",\n  \"45+18=63\": \"7+46=53\",\n  \"62+17=79\": \"35+35=70\",\n  \"61-50=11\": \"79+19=98\",\n  \"41-6=35\": \"12+77=89\",\n  \"64-25=39\": \"27-24=3\",\n  \"78-5=73\": \"58-44=14\",\n  \"38-22=16\": \"81-53=28\",\n  \"47-5=42\": \"68+31=99\",\n  \"60-24=36\": \"45+7=52\",\n  \"85-61=24\": \"19+7=26\",\n  \"77-11=66\": \"46+30=76\",\n  \"11+3=14\": \"61-41=20\",\n  \"49+26=75\": \"99-50=49\",\n  \"22+44=66\": \"97-51=46\",\n  \"91-24=67\": \"31+63=94\",\n  \"46+40=86\": \"79-26=53\",\n  \"32+20=52\": \"12+42=54\",\n  \"76-65=11\": \"20+35=55\",\n  \"48-37=11\": \"38+38=76\",\n  \"57+18=75\": \"96-1=95\",\n  \"78+16=94\": \"77-52=25\",\n  \"25-15=10\": \"74-13=61\",\n  \"44+19=63\": \"15+25=40\",\n  \"29+15=44\": \"27+64=91\",\n  \"71-53=18\": \"34+21=55\",\n  \"73-13=60\": \"58-22=36\",\n  \"47+35=82\": \"76+18=94\",\n  \"28+58=86\": \"76-7=69\",\n  \"10+24=34\": \"48+12=60\",\n  \"54+34=88\": \"37+16=53\",\n  \"11+42=53\": \"26+9=35\",\n  \"97-3=94\": \"42-29=13\",\n  \"32+19=51\": \"8+6=14\",\n  \"0+25=25\": \"85-17=68\",\n  \"36-1=35\": \"62-17=45\",\n  \"73-48=25\": \"66-39=27\",\n  \"61-56=5\": \"15+26=41\",\n  \"99-45=54\": \"98-51=47\"\n};\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document body.\");\n}\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst oldValues = table.values;\nconst newValues = oldValues.map(row =>\n  row.map(cell => {\n    const trimmed = cell.trim();\n    return Object.prototype.hasOwnProperty.call(replacements, trimmed)\n      ? replacements[trimmed]\n      : cell;\n  })\n);\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Map of old cell text -> new cell text, taken from the target diff.\n# Every data cell in the single table of this document is replaced;\n# the mapping is keyed by the *original* text so the script only\n# touches cells whose current text matches a known \"before\" value.\n$replacements = @{\n    \"68+16=84\" = \"95-44=51\"\n    \"81+2=83\" = \"14+84=98\"\n    \"54+26=80\" = \"60-35=25\"\n    \"48+20=68\" = \"4+2=6\"\n    \"24+62=86\" = \"20+43=63\"\n    \"11+13=24\" = \"27+57=84\"\n    \"7+40=47\" = \"25+21=46\"\n    \"77-1=76\" = \"81-66=15\"\n    \"59-48=11\" = \"67+23=90\"\n    \"99-81=18\" = \"20+74=94\"\n    \"4+48=52\" = \"37+61=98\"\n    \"78-50=28\" = \"30+67=97\"\n    \"56+13=69\" = \"8+73=81\"\n    \"70-6=64\" = \"94-82=12\"\n    \"70-5=65\" = \"12+32=44\"\n    \"24-6=18\" = \"35+52=87\"\n    \"37+6=43\" = \"31+61=92\"\n    \"60+32=92\" = \"32+34=66\"\n    \"34-24=10\" = \"29+61=90\"\n    \"10+76=86\" = \"58+17=75\"\n    \"62-42=20\" = \"5+26=31\"\n    \"30+12=42\" = \"68+20=88\"\n    \"32-10=22\" = \"33-0=33\"\n    \"42+49=91\" = \"34-7=27\"\n    \"23+25=48\" = \"62+11=73\"\n    \"75+17=92\" = \"37+26=63\"\n    \"41-19=22\" = \"26+47=73\"\n    \"44-31=13\" = \"63-14=49\"\n    \"76+23=99\" = \"17+76=93\"\n    \"57-39=18\" = \"93+4=97\"\n    \"58+3=61\" = \"99-22=77\"\n    \"68+5=73\" = \"42+10=52\"\n    \"7+27=34\" = \"25+52=77\"\n    \"21+12=33\" = \"43-22=21\"\n    \"32+11=43\" = \"16+53=69\"\n    \"97-90=7\" = \"71+7=78\"\n    \"44-9=35\" = \"15+45=60\"\n    \"97-28=69\" = \"49-44=5\"\n    \"99-11=88\" = \"69-39=30\"\n    \"67-13=54\" = \"57-18=39\"\n    \"23+71=94\" = \"92-71=21\"\n    \"5+9=14\" = \"55+28=83\"\n    \"38+12=50\" = \"54-17=37\"\n    \"82-81=1\" = \"66+23=89\"\n    \"58-6=52\" = \"74-65=9\"\n    \"14+15=29\" = \"42-36=6\"\n    \"2+88=90\" = \"96-92=4\"\n    \"28+43=71\" = \"56-54=2\"\n    \"69-38=31\" = \"75+10=85\"\n    \"8+51=59\" = \"11-7=4\"\n    \"93-78=15\" = \"3+48=51\"\n    \"79-52=27\" = \"35+16=51\"\n    \"86-14=72\" = \"21-8=13\"\n    \"10+82=92\" = \"96-6=90\"\n    \"58-38=20\" = \"60+20=80\"\n    \"72-41=31\" = \"55-11=44\"\n    \"57-13=44\" = \"68+23=91\"\n    \"45+2=47\" = \"63-60=3\"\n    \"76-71=5\" = \"57-35=22\"\n    \"49-46=3\" = \"61-4=57\"\n    \"17+33=50\" = \"86+9=95\"\n    \"19+12=31\" = \"65+17=82\"\n    \"45+18=63\" = \"7+46=53\"\n    \"62+17=79\" = \"35+35=70\"\n    \"61-50=11\" = \"79+19=98\"\n    \"41-6=35\" = \"12+77=89\"\n    \"64-25=39\" = \"27-24=3\"\n    \"78-5=73\" = \"58-44=14\"\n    \"38-22=16\" = \"81-53=28\"\n    \"47-5=42\" = \"68+31=99\"\n    \"60-24=36\" = \"45+7=52\"\n    \"85-61=24\" = \"19+7=26\"\n    \"77-11=66\" = \"46+30=76\"\n    \"11+3=14\" = \"61-41=20\"\n    \"49+26=75\" = \"99-50=49\"\n    \"22+44=66\" = \"97-51=46\"\n    \"91-24=67\" = \"31+63=94\"\n    \"46+40=86\" = \"79-26=53\"\n    \"32+20=52\" = \"12+42=54\"\n    \"76-65=11\" = \"20+35=55\"\n    \"48-37=11\" = \"38+38=76\"\n    \"57+18=75\" = \"96-1=95\"\n    \"78+16=94\" = \"77-52=25\"\n    \"25-15=10\" = \"74-13=61\"\n    \"44+19=63\" = \"15+25=40\"\n    \"29+15=44\" = \"27+64=91\"\n    \"71-53=18\" = \"34+21=55\"\n    \"73-13=60\" = \"58-22=36\"\n    \"47+35=82\" = \"76+18=94\"\n    \"28+58=86\" = \"76-7=69\"\n    \"10+24=34\" = \"48+12=60\"\n    \"54+34=88\" = \"37+16=53\"\n    \"11+42=53\" = \"26+9=35\"\n    \"97-3=94\" = \"42-29=13\"\n    \"32+19=51\" = \"8+6=14\"\n    \"0+25=25\" = \"85-17=68\"\n    \"36-1=35\" = \"62-17=45\"\n    \"73-48=25\" = \"66-39=27\"\n    \"61-56=5\" = \"15+26=41\"\n    \"99-45=54\" = \"98-51=47\"\n}\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $t.Cell($r, $c)\n        $raw = $cell.Range.Text\n        $key = $raw.TrimEnd([char]13, [char]7)\n        if ($replacements.ContainsKey($key)) {\n            $cell.Range.Text = $replacements[$key]\n        }\n    }\n}\n"}
